$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.120.06'
$ws.Range('E2').Value = '  +0.54%  '

$ws.Range('D3').Value = '3.336.75'
$ws.Range('E3').Value = '  +0.99%  '

$ws.Range('E4').Value = '  -0.17%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '585.67'
$ws.Range('E5').Value = '  +5.32%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '185.05'
$ws.Range('E6').Value = '  -1.11%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.999'
$ws.Range('E7').Value = '  -0.03%  '

$ws.Range('E8').Value = '  -1.45%  '

$ws.Range('E9').Value = '  -0.83%  '

$ws.Range('E10').Value = '  -0.78%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '46.91'
$ws.Range('E11').Value = '  -0.77%  '

$ws.Range('E12').Value = '  -0.42%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '657.63'
$ws.Range('E13').Value = '  +8.79%  '

$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '8.47'
$ws.Range('E14').Value = '  -2.65%  '

$ws.Range('B15').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C15').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D15').Value = '3.642.91'
$ws.Range('E15').Value = '  -4.92%  '

$ws.Range('D16').Value = '66.289.17'
$ws.Range('E16').Value = '  +0.80%  '

$ws.Range('E17').Value = '  -0.16%  '

$ws.Range('D18').Value = '3.336.98'
$ws.Range('E18').Value = '  +0.95%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '17.85'
$ws.Range('E19').Value = '  -0.76%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.09'
$ws.Range('E20').Value = '  +0.13%  '

$ws.Range('E21').Value = '  -1.39%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '17.61'
$ws.Range('E22').Value = '  -5.01%  '

$ws.Range('E23').Value = '  -0.55%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '100.21'
$ws.Range('E24').Value = '  +0.01%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '4.00'
$ws.Range('E25').Value = '  +1.46%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.78'
$ws.Range('E26').Value = '  +0.13%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.51'
$ws.Range('E27').Value = '  -0.68%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '32.07'
$ws.Range('E28').Value = '  +5.85%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.50'
$ws.Range('E29').Value = '  -2.11%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '6.73'
$ws.Range('E30').Value = '  -0.20%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '604.87'
$ws.Range('E31').Value = '  +2.68%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.88'
$ws.Range('E32').Value = '  +1.75%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '11.07'
$ws.Range('E33').Value = '  -0.35%  '

$ws.Range('D34').Value = '3.882.62'
$ws.Range('E34').Value = '  +4.85%  '

$ws.Range('E35').Value = '  +0.22%  '

$ws.Range('E36').Value = '  +0.02%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '56.32'
$ws.Range('E37').Value = '  -1.18%  '

$ws.Range('E38').Value = '  -1.37%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.71'
$ws.Range('E39').Value = '  +1.55%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '32.87'
$ws.Range('E40').Value = '  -3.35%  '

$ws.Range('D41').Value = '0.0₃0698'
$ws.Range('E41').Value = '  -3.29%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.18'
$ws.Range('E42').Value = '  -2.69%  '

$ws.Range('E43').Value = '  +1.25%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.337'
$ws.Range('E44').Value = '  -1.25%  '

$ws.Range('E45').Value = '  -1.14%  '

$ws.Range('E46').Value = '  -1.41%  '

$ws.Range('E47').Value = '  +0.40%  '

$ws.Range('E48').Value = '  -1.66%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.83'
$ws.Range('E49').Value = '  -18.11%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.33'
$ws.Range('E50').Value = '  +6.42%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '129.29'
$ws.Range('E51').Value = '  +4.38%  '
